$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- plain value edits (no hyperlink relationship) ---
$ws.Range("O34").Value = "https://wallpapercave.com/wp/wp7480220.jpg"
$ws.Range("B44").Value = "Crystal Dynamics, Nixxes Software, Feral Interactive, Buzz Monkey Software, Santa Cruz Games"
$ws.Range("C44").Value = 2008
$ws.Range("B45").Value = "Crystal Dynamics, Eidos-Montréal, Nixxes"
$ws.Range("C45").Value = 2015
$ws.Range("B46").Value = "Crystal Dynamics, Eidos Montreal, Nixxes Software"
$ws.Range("C46").Value = 2018
$ws.Range("B65").Value = "Dice"
$ws.Range("C99").Value = 2008
$ws.Range("C100").Value = 2013
$ws.Range("B101").Value = "Crytek"
$ws.Range("C116").Value = 2015
$ws.Range("O116").Value = "https://i.pinimg.com/originals/ff/7d/d6/ff7dd6249fd6745653fc37a8d014aa01.jpg"
$ws.Range("B120").Value = "Blizzard Entertainment, Blizzard North"
$ws.Range("C120").Value = 1996
$ws.Range("B121").Value = "Blizzard Entertainment, Vicarious Visions, Blizzard North"
$ws.Range("C121").Value = 2000
$ws.Range("B122").Value = "Blizzard Entertainment"
$ws.Range("C122").Value = 2012
$ws.Range("D143").Value = "Co-op, puzzle"

# --- hyperlinks: bordered-style group first (fixes style slot 32) ---
$ws.Hyperlinks.Add($ws.Range("O35"), "https://images8.alphacoders.com/107/thumb-1920-1078901.jpg") | Out-Null

# --- seed the non-bordered style (slot 33) ---
$ws.Hyperlinks.Add($ws.Range("O89"), "https://wallpaperaccess.com/full/672018.jpg") | Out-Null

# --- remaining bordered-group hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("O36"), "https://www.well-played.com.au/wp-content/uploads/2020/05/m2logo.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O37"), "https://images.hdqwalls.com/wallpapers/mafia-iii-definitive-edition-be.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O38"), "https://images.hdqwalls.com/wallpapers/call-of-duty-modern-warfare-remastered-lu.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O39"), "https://images.hdqwalls.com/wallpapers/call-of-duty-modern-warfare-2-campaign-remastered-4k-vw.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O40"), "https://images2.alphacoders.com/939/thumb-1920-939737.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O41"), "https://static.wikia.nocookie.net/residentevil/images/1/11/About_bg_re3.jpg/revision/latest?cb=20210226115827&path-prefix=ru") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O42"), "https://i.pinimg.com/originals/e7/31/51/e7315109ffc6599e2a5372fc95b24ba4.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O44"), "https://images7.alphacoders.com/330/thumb-1920-330535.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O45"), "https://images5.alphacoders.com/724/724192.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O46"), "https://i.ytimg.com/vi/TgAyj-mtDg4/maxresdefault.jpg") | Out-Null

# --- remaining non-bordered-group hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("O91"), "https://i.pinimg.com/originals/b4/50/ca/b450ca40de0bd9849cbc0fdd67163883.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O92"), "https://images2.alphacoders.com/509/thumb-1920-509156.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O93"), "https://wallpapercave.com/wp/wp2867631.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O94"), "https://images7.alphacoders.com/605/thumb-1920-605394.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O99"), "https://wallpaperaccess.com/full/1482238.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O100"), "https://wallpaperaccess.com/full/957573.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O101"), "https://wallpaperaccess.com/full/3270333.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O102"), "https://images.alphacoders.com/148/thumb-1920-148221.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O103"), "https://i.pinimg.com/originals/0d/7f/44/0d7f44c79af0b05f2b8ae51cd812a04c.png") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O104"), "https://images2.alphacoders.com/664/664487.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O105"), "https://ixbt.online/gametech/games/2021/02/01/bf5ygoPgEzF3GVVBMitRZdB1uIMxkeUM9qUWj7lG.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O106"), "https://i.pinimg.com/originals/28/36/33/2836330c4eb24377b5c4794669a2ec41.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O107"), "https://i.playground.ru/p/FcZcnV95D0IyFzSNuuLPnw.jpeg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O111"), "https://images7.alphacoders.com/800/thumb-1920-800857.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O114"), "https://wallpaperaccess.com/full/653677.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O115"), "https://images3.alphacoders.com/828/thumb-1920-828135.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O120"), "https://cdn.mos.cms.futurecdn.net/G5sChiPXCGDC855Ao78HxR.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O121"), "https://ixbt.online/gametech/covers/2021/04/24/HwYAEJMXBGhA2sk1o8WipeuxPPDOsY38o2OUG5S6.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O122"), "https://wallpapercave.com/wp/6IyEN0s.jpg") | Out-Null

# --- O90: gets a hyperlink but keeps its original bordered style (s=2) ---
$ws.Hyperlinks.Add($ws.Range("O90"), "https://i.pinimg.com/originals/31/01/f6/3101f69af8406ff979d8b4c1032badbe.jpg") | Out-Null
$ws.Range("N90").Copy() | Out-Null
$ws.Range("O90").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

Write-Output "edits applied"
